# Applies the "heavy loadtest data corrected" edit described by the commit:
#   fix: heavy loadtest data corrected according to confluence Performance Test Data.xls
#
# Changes:
#  - heavy!E6  100   -> 0
#  - heavy!G6  0.01  -> 0
#  - heavy!H6  =100+K6 (formula) -> 0 (plain value)
#  - heavy!M6  0.01  -> 0
#  - heavy!G7  0.24  -> 0.25
#  - heavy!M7  0.24  -> 0.25
#  - makes "heavy" sheet the active/selected sheet & tab (it was "light" before),
#    with the active selection on heavy set to M8.

$wb = $excel.ActiveWorkbook
$heavy = $wb.Worksheets.Item("heavy")

$heavy.Range("E6").Value = 0
$heavy.Range("G6").Value = 0
$heavy.Range("H6").Value = 0
$heavy.Range("M6").Value = 0
$heavy.Range("G7").Value = 0.25
$heavy.Range("M7").Value = 0.25

# Recalculate all formulas so dependent cells (F6, I6, J6, I5, F7, I7, and the
# row-2 summary formulas) reflect the edits above.
$wb.Application.Calculate()

# Switch the active sheet/tab from "light" to "heavy" and set the selection.
$heavy.Activate()
$heavy.Range("M8").Select()
